# Update TODO list: 2021.05.17.
# - Replace the placeholder letters (a..j) with task codes (T10010..T10090)
# - Drop the last row (row 10), shrinking the used range to A1:B9
# - Center the task-code cells vertically (A2:B9)
# - Move the active selection to H5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused last row first so the remaining data shifts
# up cleanly and the sheet dimension becomes A1:B9.
$ws.Rows(10).Delete()

$values = @("T10010", "T10020", "T10030", "T10040", "T10050", "T10060", "T10070", "T10080", "T10090")

for ($i = 0; $i -lt ($values.Length - 1); $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i + 1]
}

# Vertically center the task-code cells (new style: xf with applyAlignment vertical=center)
$ws.Range("A2:B9").VerticalAlignment = -4108

# Update the active selection
$ws.Range("H5").Select()
